$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from its current location (end of paragraph 3);
#    it will be re-inserted at the top of paragraph 1 below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Paragraph 1 ("Arpeggi"): switch from manual bold/size run formatting to the
#    built-in "Titolo1" (Heading 1) paragraph style, and plant the _GoBack
#    bookmark right at the top of the paragraph.
$p1 = $d.Paragraphs(1)
$xml1 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Titolo1"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Arpeggi</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p1.Range.InsertXML($xml1)

# 3. Paragraph 2: rewrite the body text (drop "che nelle scale, " and turn
#    "favore" into "favorire"), splitting the sentence across three runs.
$p2 = $d.Paragraphs(2)
$xml2 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:xml="http://www.w3.org/XML/1998/namespace"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Il meccanismo del pollice richiede, negli arpeggi, un orientamento della mano al quanto pi&#249; accurato e sar&#224; anche necessario ricorrere ogni volta a una leggera rotazione del polso per favor</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ire</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> il passaggio del pollice. I muscoli della mano destra dovranno essere rilassati come nelle scale, seppur sia pi&#249; difficile, data la maggior tensione muscolare richiesta dalla particolare configurazione degli arpeggi.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($xml2)
